$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Header row (row 1) restyle: bold black text on a gold highlight fill.
# ---------------------------------------------------------------------------
$header = $ws.Range("A1:M1")
$header.Font.Bold = $true
$header.Font.Color = 0
$header.Interior.Color = 10086143
$header.Interior.PatternColor = 0

# ---------------------------------------------------------------------------
# Existing order row (row 2): OrderID format changed, and the whole row now
# uses the data font (fontId 19) with a distinct border on the last column.
# ---------------------------------------------------------------------------
$ws.Range("A2").NumberFormat = "m/d/yy h:mm"
$ws.Range("A2").Value = 45436

$ws.Range("B2").Value = "001BC"
$ws.Range("C2").Value = "SIPL5316"
$ws.Range("D2").Value = "SIPL5688"
$ws.Range("E2").Value = "SIPL0102"
$ws.Range("F2").Value = "SIPL0103"
$ws.Range("G2").Value = "BC Law Firm"
$ws.Range("H2").Value = "Title"
$ws.Range("I2").Value = "Search & Typing"
$ws.Range("J2").Value = "Full Search"
$ws.Range("K2").Value = "FL"
$ws.Range("L2").Value = "Clay"
$ws.Range("M2").Value = "WIP"

$ws.Range("A2:M2").Font.Color = 0

# ---------------------------------------------------------------------------
# New order rows 3 and 4.
# ---------------------------------------------------------------------------
$ws.Range("A3").NumberFormat = "m/d/yy h:mm"
$ws.Range("A3").Value = 45439
$ws.Range("B3").Value = "002BC"
$ws.Range("C3").Value = "SIPL5316"
$ws.Range("D3").Value = "SIPL5688"
$ws.Range("E3").Value = "SIPL0102"
$ws.Range("F3").Value = "SIPL0103"
$ws.Range("G3").Value = "BC Law Firm"
$ws.Range("H3").Value = "Title"
$ws.Range("I3").Value = "Search & Typing"
$ws.Range("J3").Value = "Update Search"
$ws.Range("K3").Value = "FL"
$ws.Range("L3").Value = "Clay"
$ws.Range("M3").Value = "WIP"

$ws.Range("A4").NumberFormat = "m/d/yy h:mm"
$ws.Range("A4").Value = 45439
$ws.Range("B4").Value = "003BC"
$ws.Range("C4").Value = "SIPL5316"
$ws.Range("D4").Value = "SIPL5688"
$ws.Range("E4").Value = "SIPL0102"
$ws.Range("F4").Value = "SIPL0103"
$ws.Range("G4").Value = "BC Law Firm"
$ws.Range("H4").Value = "Title"
$ws.Range("I4").Value = "Search & Typing"
$ws.Range("J4").Value = "Current Owner Search"
$ws.Range("K4").Value = "FL"
$ws.Range("L4").Value = "Clay"
$ws.Range("M4").Value = "WIP"

$ws.Range("A3:M4").Font.Color = 0

# ---------------------------------------------------------------------------
# Borders: last column of each data row (D) gets the alternate border.
# ---------------------------------------------------------------------------
$ws.Range("D2:D4").Borders.LineStyle = 1

# Selection, matching the saved workbook's cursor position.
$ws.Range("G4").Select()
